# Add a new "2022-Q4" worksheet right after "总计", shifting the existing
# quarter sheets (2021-Q2, 2021-Q1, 2020-Q4) one position to the right, and
# update the "总计" (summary) sheet with the new quarter's totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet: insert a new row for 2022-Q4 at the
#    top of the data, pushing the existing rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Shift the B:D data down one row (bottom-up so we don't clobber data
# before it is copied). This also carries the correct cell formatting.
$total.Range("B4:D4").Copy($total.Range("B5:D5"))
$total.Range("B3:D3").Copy($total.Range("B4:D4"))
$total.Range("B2:D2").Copy($total.Range("B3:D3"))

# Extend the styled index column (A) down to the new row 5.
$total.Range("A4").Copy($total.Range("A5"))

# Re-number the sequential index column.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Fill in the new 2022-Q4 totals row.
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 9
$total.Range("D2").Value = 0.34

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Copy header + row formatting from the existing "2021-Q2" sheet so the
# new sheet matches the look (bold header, bordered index column).
$srcFormat = $wb.Worksheets.Item("2021-Q2")
$srcFormat.Range("B1:H2").Copy($q4.Range("B1:H2"))

# Header row.
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Replicate the styled index-column cell (bold + bordered) down to rows 2-10.
$srcFormat.Range("A2").Copy($q4.Range("A2:A10"))

# Fund holding rows.
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'005825"
$q4.Range("C2").Value = "申万菱信智能驱动股票A"
$q4.Range("D2").Value = "'6.22"
$q4.Range("E2").Value = "'84.52"
$q4.Range("F2").Value = "'2.72"
$q4.Range("G2").Value = "'0.1692"
$q4.Range("H2").Value = 8

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'015159"
$q4.Range("C3").Value = "申万菱信智能驱动股票C"
$q4.Range("D3").Value = "'2.08"
$q4.Range("E3").Value = "'84.52"
$q4.Range("F3").Value = "'2.72"
$q4.Range("G3").Value = "'0.0566"
$q4.Range("H3").Value = 8

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'005310"
$q4.Range("C4").Value = "广发电子信息传媒产业精选股票A"
$q4.Range("D4").Value = "'1.72"
$q4.Range("E4").Value = "'88.90"
$q4.Range("F4").Value = "'2.93"
$q4.Range("G4").Value = "'0.0504"
$q4.Range("H4").Value = 7

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'009048"
$q4.Range("C5").Value = "浦银安盛科技创新优选三年封闭运作灵活配置混合"
$q4.Range("D5").Value = "'2.39"
$q4.Range("E5").Value = "'35.76"
$q4.Range("F5").Value = "'1.07"
$q4.Range("G5").Value = "'0.0256"
$q4.Range("H5").Value = 8

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'015919"
$q4.Range("C6").Value = "申万菱信专精特新主题混合A"
$q4.Range("D6").Value = "'0.39"
$q4.Range("E6").Value = "'48.55"
$q4.Range("F6").Value = "'3.85"
$q4.Range("G6").Value = "'0.0150"
$q4.Range("H6").Value = 7

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'519127"
$q4.Range("C7").Value = "浦银安盛盛世精选灵活配置混合A"
$q4.Range("D7").Value = "'1.17"
$q4.Range("E7").Value = "'23.43"
$q4.Range("F7").Value = "'0.73"
$q4.Range("G7").Value = "'0.0085"
$q4.Range("H7").Value = 8

$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "'010236"
$q4.Range("C8").Value = "广发电子信息传媒产业精选股票C"
$q4.Range("D8").Value = "'0.24"
$q4.Range("E8").Value = "'88.90"
$q4.Range("F8").Value = "'2.93"
$q4.Range("G8").Value = "'0.0070"
$q4.Range("H8").Value = 7

$q4.Range("A9").Value = 7
$q4.Range("B9").Value = "'519177"
$q4.Range("C9").Value = "浦银安盛盛世精选灵活配置混合C"
$q4.Range("D9").Value = "'0.76"
$q4.Range("E9").Value = "'23.43"
$q4.Range("F9").Value = "'0.73"
$q4.Range("G9").Value = "'0.0055"
$q4.Range("H9").Value = 8

$q4.Range("A10").Value = 8
$q4.Range("B10").Value = "'015920"
$q4.Range("C10").Value = "申万菱信专精特新主题混合C"
$q4.Range("D10").Value = "'0.02"
$q4.Range("E10").Value = "'48.55"
$q4.Range("F10").Value = "'3.85"
$q4.Range("G10").Value = "'0.0008"
$q4.Range("H10").Value = 7

# Restore the originally-active "2020-Q4" tab (adding a sheet makes the
# newly inserted sheet active by default).
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "Done: added 2022-Q4 sheet and updated totals."
